$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2548.4167
$ws.Range("I40").Value = 5650.5
$ws.Range("J40").Value = 1928
$ws.Range("K40").Value = 5650.5
$ws.Range("L40").Value = 1928
$ws.Range("M40").Value = -5475.5
$ws.Range("N40").Value = -2278
# Row 76
$ws.Range("H76").Value = 5574.68
$ws.Range("I76").Value = 4542.778
$ws.Range("J76").Value = 6155.125
$ws.Range("K76").Value = 4542.778
$ws.Range("L76").Value = 6155.125
$ws.Range("M76").Value = -4227.778
$ws.Range("N76").Value = -6785.125
# Row 79
$ws.Range("H79").Value = 5574.68
$ws.Range("I79").Value = 4542.778
$ws.Range("J79").Value = 6155.125
$ws.Range("K79").Value = 4542.778
$ws.Range("L79").Value = 6155.125
$ws.Range("M79").Value = -3450.778
$ws.Range("N79").Value = -8339.125
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
# Row 132
$ws.Range("H132").Value = 8776922
$ws.Range("I132").Value = 10422132
$ws.Range("K132").Value = 31266396
$ws.Range("M132").Value = -31263866
# Row 138
$ws.Range("H138").Value = 545352.4
$ws.Range("I138").Value = 1966.5
$ws.Range("J138").Value = 641244
$ws.Range("K138").Value = 5899.5
$ws.Range("L138").Value = 1923732
$ws.Range("M138").Value = -759.5
$ws.Range("N138").Value = -1934012
# Row 141
$ws.Range("H141").Value = 300
$ws.Range("I141").Value = 300
$ws.Range("K141").Value = 900
$ws.Range("M141").Value = 4280

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2856.817
$ws.Range("I32").Value = 2677.4688
$ws.Range("K32").Value = 2677.4688
$ws.Range("M32").Value = -2390.4688
# Row 74
$ws.Range("H74").Value = 1026.3334
$ws.Range("I74").Value = 492.77274
$ws.Range("J74").Value = 2093.4546
$ws.Range("K74").Value = 492.77274
$ws.Range("L74").Value = 2093.4546
$ws.Range("M74").Value = 381.22726
$ws.Range("N74").Value = -3841.4546
# Row 77
$ws.Range("H77").Value = 1026.3334
$ws.Range("I77").Value = 492.77274
$ws.Range("J77").Value = 2093.4546
$ws.Range("K77").Value = 2463.8637
$ws.Range("L77").Value = 10467.273
$ws.Range("M77").Value = 1904.1363
$ws.Range("N77").Value = -19203.273
# Row 97
$ws.Range("H97").Value = 307.53845
$ws.Range("I97").Value = 318.16666
$ws.Range("J97").Value = 180
$ws.Range("K97").Value = 318.16666
$ws.Range("L97").Value = 180
$ws.Range("M97").Value = 177.83334
$ws.Range("N97").Value = -1172
# Row 110
$ws.Range("H110").Value = 1035.1666
$ws.Range("I110").Value = 807.2778
$ws.Range("K110").Value = 807.2778
$ws.Range("M110").Value = 1237.7222

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5608.577
$ws.Range("I134").Value = 1037.409
$ws.Range("K134").Value = 3112.227
$ws.Range("M134").Value = -577.2270000000003
# Row 138
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1911.1538
$ws.Range("I31").Value = 1585.3334
$ws.Range("J31").Value = 2355.4546
$ws.Range("K31").Value = 1585.3334
$ws.Range("L31").Value = 2355.4546
$ws.Range("M31").Value = -1290.3334
$ws.Range("N31").Value = -2945.4546
# Row 34
$ws.Range("H34").Value = 1911.1538
$ws.Range("I34").Value = 1585.3334
$ws.Range("J34").Value = 2355.4546
$ws.Range("K34").Value = 1585.3334
$ws.Range("L34").Value = 2355.4546
$ws.Range("M34").Value = -1383.3334
$ws.Range("N34").Value = -2759.4546
# Row 105
$ws.Range("H105").Value = 793.1667
$ws.Range("I105").Value = 751.8
$ws.Range("K105").Value = 751.8
$ws.Range("M105").Value = 995.2
# Row 114
$ws.Range("H114").Value = 24775.555
$ws.Range("J114").Value = 24775.555
$ws.Range("L114").Value = 24775.555
$ws.Range("N114").Value = -33453.555
# Row 132
$ws.Range("H132").Value = 6906.9546
$ws.Range("I132").Value = 9480.846
$ws.Range("J132").Value = 3189.111
$ws.Range("K132").Value = 28442.538
$ws.Range("L132").Value = 9567.332999999999
$ws.Range("M132").Value = -25912.538
$ws.Range("N132").Value = -14627.333
# Row 134
$ws.Range("H134").Value = 1896.8572
$ws.Range("I134").Value = 1847.6364
$ws.Range("K134").Value = 5542.9092
$ws.Range("M134").Value = -3007.9092

$ws = $wb.Worksheets.Item("CUL")
# Row 136
$ws.Range("H136").Value = 1406
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 137
$ws.Range("H137").Value = 1906.75
$ws.Range("I137").Value = 1168
$ws.Range("J137").Value = 10033
$ws.Range("K137").Value = 3504
$ws.Range("L137").Value = 30099
$ws.Range("M137").Value = 1596
$ws.Range("N137").Value = -40299
# Row 138
$ws.Range("H138").Value = 3211.125
$ws.Range("I138").Value = 3348.1667
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 10044.5001
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = -4904.500100000001
$ws.Range("N138").Value = -18680
# Row 139
$ws.Range("H139").Value = 1575.7435
$ws.Range("I139").Value = 1629.9584
$ws.Range("K139").Value = 4889.8752
$ws.Range("M139").Value = 250.1247999999996
# Row 140
$ws.Range("H140").Value = 33357.742
$ws.Range("I140").Value = 43876.117
$ws.Range("K140").Value = 131628.351
$ws.Range("M140").Value = -126448.351

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4001.0908
$ws.Range("I80").Value = 3576.5
$ws.Range("J80").Value = 5133.3335
$ws.Range("K80").Value = 3576.5
$ws.Range("L80").Value = 5133.3335
$ws.Range("M80").Value = -2578.5
$ws.Range("N80").Value = -7129.3335
# Row 83
$ws.Range("H83").Value = 4001.0908
$ws.Range("I83").Value = 3576.5
$ws.Range("J83").Value = 5133.3335
$ws.Range("K83").Value = 17882.5
$ws.Range("L83").Value = 25666.6675
$ws.Range("M83").Value = -12890.5
$ws.Range("N83").Value = -35650.6675
# Row 102
$ws.Range("H102").Value = 2418.425
$ws.Range("I102").Value = 2329.6667
$ws.Range("K102").Value = 2329.6667
$ws.Range("M102").Value = -707.6667000000002
# Row 107
$ws.Range("H107").Value = 755.2857
$ws.Range("J107").Value = 476.25
$ws.Range("L107").Value = 476.25
$ws.Range("N107").Value = -4316.25
# Row 122
$ws.Range("H122").Value = 1337.2122
$ws.Range("I122").Value = 1459.2084
$ws.Range("J122").Value = 1011.8889
$ws.Range("K122").Value = 4377.6252
$ws.Range("L122").Value = 3035.6667
$ws.Range("M122").Value = -1927.6252
$ws.Range("N122").Value = -7935.6667
# Row 126
$ws.Range("H126").Value = 1991.6666
$ws.Range("I126").Value = 1657.9286
$ws.Range("K126").Value = 4973.7858
$ws.Range("M126").Value = -2503.7858
# Row 141
$ws.Range("H141").Value = 46833.332
$ws.Range("J141").Value = 46833.332
$ws.Range("L141").Value = 46833.332
$ws.Range("N141").Value = -57193.332

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 90910470
$ws.Range("I126").Value = 100001470
$ws.Range("J126").Value = 495
$ws.Range("K126").Value = 300004410
$ws.Range("L126").Value = 1485
$ws.Range("M126").Value = -300001940
$ws.Range("N126").Value = -6425
# Row 128
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
# Row 132
$ws.Range("H132").Value = 1737.2051
$ws.Range("I132").Value = 1238.0385
$ws.Range("K132").Value = 3714.1155
$ws.Range("M132").Value = -1184.1155
# Row 136
$ws.Range("H136").Value = 381.21213
$ws.Range("I136").Value = 287.96295
$ws.Range("J136").Value = 800.8333
$ws.Range("K136").Value = 863.8888499999999
$ws.Range("L136").Value = 2402.4999
$ws.Range("M136").Value = 1686.11115
$ws.Range("N136").Value = -7502.4999
# Row 137
$ws.Range("H137").Value = 30000
$ws.Range("J137").Value = 30000
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200
